$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: date serial 45310 -> 45311 (2024-01-19 -> 2024-01-20)
$ws.Range("A1").Value = 45311

# D29: 364.992 -> 185.28
$ws.Range("D29").Value = 185.28

# D30: 514.29 -> 261.067
$ws.Range("D30").Value = 261.067
